# Refresh the sliding-window inference results (window 4): IPC PO predictions
# in column C for rows 2-51, along with their dependent DELTA (col D) and
# DELTA^2 (col E) values, and the TOTAL (row 52) / MSE (row 53) summary cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, IPC PO (col C), DELTA (col D), DELTA^2 (col E)
$data = @(
    ,@(2, 29.93040638646515, 0.5604063864651536, 0.3140553179909311)
    ,@(3, 29.32506259951632, -0.2149374004836773, 0.04619808612668067)
    ,@(4, 29.77535200065974, 0.2253520006597398, 0.05078352420134737)
    ,@(5, 29.6515993714725, -0.09840062852750009, 0.009682683694607064)
    ,@(6, 30.2031358527179, 0.3631358527178961, 0.1318676475291535)
    ,@(7, 30.08469355855535, 0.2746935585553487, 0.07545655111180079)
    ,@(8, 29.83023722208232, -0.08976277791768439, 0.008057356299499528)
    ,@(9, 29.75417165550755, -0.2258283444924487, 0.05099844117620009)
    ,@(10, 29.76025260776232, -0.279747392237681, 0.07825860346378292)
    ,@(11, 29.83803663878533, -0.371963361214668, 0.1383567420861136)
    ,@(12, 30.01462314165534, -0.2053768583446605, 0.04217965394352273)
    ,@(13, 30.31141795601553, -0.06858204398446688, 0.004703496757087349)
    ,@(14, 30.57077553773555, 0.1307755377355484, 0.01710224127002184)
    ,@(15, 30.32483362114042, -0.1551663788595832, 0.0240766051283957)
    ,@(16, 30.35613061192455, -0.3338693880754562, 0.1114687682938796)
    ,@(17, 30.87410749143421, 0.1241074914342057, 0.01540266943009145)
    ,@(18, 30.7632536985891, -0.1767463014108976, 0.03123925506243188)
    ,@(19, 30.89930378231916, -0.05069621768083721, 0.002570106487142831)
    ,@(20, 31.13777421556491, 0.1177742155649142, 0.01387076585193087)
    ,@(21, 31.28926585984339, 0.1692658598433923, 0.02865093130852291)
    ,@(22, 31.19905793677455, -0.08094206322545361, 0.00655161759919333)
    ,@(23, 31.11989792002061, -0.2601020799793901, 0.06765309200960502)
    ,@(24, 31.32230849825455, -0.2576915017454482, 0.06640491007182434)
    ,@(25, 31.78016681862083, 0.130166818620836, 0.01694340066986962)
    ,@(26, 32.52623230257632, 0.646232302576319, 0.4176161888930912)
    ,@(27, 32.57684602740945, 0.296846027409444, 0.08811756398876837)
    ,@(28, 32.91880862076341, 0.4688086207634043, 0.2197815229020855)
    ,@(29, 33.19745864027098, 0.3474586402709789, 0.1207275066989576)
    ,@(30, 33.25509844645786, 0.3550984464578661, 0.12609490667679)
    ,@(31, 33.33737814353816, 0.2373781435381588, 0.05634838302962272)
    ,@(32, 33.74769882878559, 0.3476988287855889, 0.1208944755388703)
    ,@(33, 33.74683683692324, 0.04683683692324081, 0.002193689292974254)
    ,@(34, 33.81645724358367, -0.2835427564163311, 0.08039649471617086)
    ,@(35, 34.20987557560572, -0.190124424394277, 0.03614729675125516)
    ,@(36, 34.42951112912753, -0.4704888708724653, 0.2213597776148473)
    ,@(37, 35.46831273599246, 0.1683127359924583, 0.02832917709726698)
    ,@(38, 35.73577698497842, 0.03577698497841908, 0.001279992654146024)
    ,@(39, 36.01183120736154, -0.2881687926384586, 0.08304125305070696)
    ,@(40, 36.63353645952962, -0.1664635404703745, 0.027710110305932)
    ,@(41, 36.72230890903101, -0.5776910909689903, 0.3337269965849422)
    ,@(42, 37.82487224815392, -0.07512775184607534, 0.005644179097445477)
    ,@(43, 38.59315908521873, 0.09315908521872984, 0.00867861515879057)
    ,@(44, 39.04101527828282, 0.141015278282822, 0.01988530870918173)
    ,@(45, 39.44285168523187, 0.04285168523186655, 0.00183626692721097)
    ,@(46, 39.70063914821225, -0.1993608517877519, 0.03974474922553797)
    ,@(47, 40.02057981949194, -0.07942018050805899, 0.006307565071932674)
    ,@(48, 40.24804694129041, -0.35195305870959, 0.1238709555350361)
    ,@(49, 40.5721519782377, -0.3278480217622999, 0.1074843253734535)
    ,@(50, 41.70425954152643, 0.5042595415264302, 0.2542776852204455)
    ,@(51, 41.83056249333569, 0.330562493335691, 0.1092715620003088)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]   # Column C: IPC PO
    $ws.Cells.Item($r, 4).Value = $row[2]   # Column D: DELTA
    $ws.Cells.Item($r, 5).Value = $row[3]   # Column E: DELTA^2
}

$ws.Cells.Item(52, 3).Value = 0.277971294333927   # C52: TOTAL of DELTA
$ws.Cells.Item(52, 5).Value = 3.993299015679407   # E52: TOTAL of DELTA^2
$ws.Cells.Item(53, 5).Value = 0.07986598031358813   # E53: MSE
